$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 / G1 hold month-name labels ("September 2024" / "October 2024").
# Assigning a plain string like "October 2024" via .Value gets
# auto-recognized by Excel as a date and converted to a date serial plus a
# new cell style. Forcing the cell to Text format before the assignment
# keeps it as literal text; clearing the format afterwards drops the
# leftover style index back to the sheet's default (General) style while
# preserving the text that was already committed to the cell.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "October 2024"
$ws.Range("A1").ClearFormats()

$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "November 2024"
$ws.Range("G1").ClearFormats()

# Row 2 data values shift over (new month appended, oldest dropped).
$ws.Range("A2").Value = 0.155
$ws.Range("B2").Value = 0.64
$ws.Range("C2").Value = -0.086
$ws.Range("D2").Value = -0.025
$ws.Range("E2").Value = 0.025
$ws.Range("F2").Value = -0.172
$ws.Range("G2").Value = 0.536
